$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the YouTube link in D11 with a hyperlink
$ws.Hyperlinks.Add($ws.Range("D11"), "https://youtu.be/U3ydTsRwxok", "", "", "https://youtu.be/U3ydTsRwxok ")
$ws.Range("D11").Style = "Collegamento ipertestuale"

# Update the term text in B11
$ws.Range("B11").Value = "Spontaneous intrahepatic porto-systemic shunt"

# Update selection to B11
$ws.Range("B11").Select()
